# Update the cryptos price/volume table (Price = column D, Volume(1h) = column E).
# For D-column values that parse as plain numbers, force the cell to Text
# ("@" number format) before assigning so the literal numeric-looking string
# (e.g. "240.64", "0.000008176", "1.0000") is preserved verbatim instead of
# being auto-coerced into a number; then restore the "Normal" style so no
# stray cell formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.066.62"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.831.26"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6843"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3012"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07454"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07662"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("D12").Value = "1.846.36"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.062"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6818"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.147"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.28%  "
$ws.Range("D17").Value = "29.077.26"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008176"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").Value = "2.083.92"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.97%  "
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.743"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.512"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.291"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.149"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.196"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05163"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7673"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.845"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.677"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("D38").Value = "1.307.74"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.726"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9311"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.813"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("D47").Value = "1.985.12"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5189"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.563"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.771"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("E51").Value = "  +0.88%  "
